$p = $ppt.ActivePresentation

# --- Slide 1: title-page subtitle / author credit line ---
# "The Official strategy guide Written by :Anne Mole"
#   -> "The Official strategy guide Written by: Anmol"
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item("Subtitle 2")
$tr1 = $sh1.TextFrame.TextRange
# Re-write via a full-length Characters() range (instead of TextRange.Text
# directly) so the shape's second, empty trailing paragraph is preserved
# and the result stays a single run.
$full1 = $tr1.Characters(1, $tr1.Length)
$full1.Text = "The Official strategy guide Written by: Anmol"

# --- Slide 3: prologue textbox ---
# Collapses the "Prince <SkullCuddle>, the second..." run split into a
# single run reading "Prince Skull, the second...".
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item("TextBox 5")
$tr3 = $sh3.TextFrame.TextRange
$full3 = $tr3.Characters(1, $tr3.Length)
$full3.Text = "For generations they ruled the land with grace and dignity. Prince Skull, the second son, upon hearing he would not become king made a deal with a Witch to turn all his subjects into monsters. In return, he would become powerful far beyond any human’s comprehension. "

# --- Slide 5: Goblin description textbox, "cure ordeal" -> "curse ordeal" ---
# Leave the leading "Goblin" run (bold) untouched; only rewrite the second
# run's text (everything after it) so it stays a single run, matching the
# diff which only edits that run's <a:t>.
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item("TextBox 3")
$tr5 = $sh5.TextFrame.TextRange
$goblinLen = 6
$rest5 = $tr5.Characters($goblinLen + 1, $tr5.Length - $goblinLen)
$rest5.Text = ": These monsters by far got the worst end of the whole curse ordeal, being given a weakened body and a panicked mind. The only thing they have left is their ability to spit rocks."
